$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the photo URL for MARCELO QUEIROGA (row 24) to the new
# upload.wikimedia.org direct image link.
$ws.Range("B24").Value = "https://upload.wikimedia.org/wikipedia/commons/f/f6/CAS_-_Comiss%C3%A3o_de_Assuntos_Sociais_%2835941643904%29_%28cropped%29.jpg"

# Match the author's final view state: scrolled down with B24 selected.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B24").Select()
